$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.010.26'
$ws.Range("E2").Value = '  +1.62%  '
$ws.Range("D3").Value = '3.426.74'
$ws.Range("E3").Value = '  +1.25%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.00'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.88'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.39%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.58'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.65%  '
$ws.Range("E10").Value = '  +0.85%  '
$ws.Range("E11").Value = '  -0.23%  '
$ws.Range("D12").Value = '4.010.51'
$ws.Range("E13").Value = '  -0.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.43'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.03%  '
$ws.Range("D15").Value = '3.426.74'
$ws.Range("E15").Value = '  +1.47%  '
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("D17").Value = '62.055.98'
$ws.Range("E17").Value = '  +1.56%  '
$ws.Range("E18").Value = '  +0.95%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.96'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.46%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.20'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.85%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '392.11'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '74.57'
$ws.Range("D22").Style = "Normal"
$ws.Range("E23").Value = '  +0.29%  '
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("E25").Value = '  +0.36%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.191'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.19%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.50'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.97%  '
$ws.Range("E28").Value = '  -0.02%  '
$ws.Range("E29").Value = '  +0.66%  '
$ws.Range("E30").Value = '  +0.79%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.41'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.82%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.60'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.61%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.31'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.75%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.99'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.44%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '168.03'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.24%  '
$ws.Range("D37").Value = '3.459.73'
$ws.Range("E38").Value = '  +1.59%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '28.44'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.39%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0757'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.19%  '
$ws.Range("E41").Value = '  +1.32%  '
$ws.Range("E42").Value = '  +1.75%  '
$ws.Range("E43").Value = '  +1.51%  '
$ws.Range("E44").Value = '  +4.66%  '
$ws.Range("D45").Value = '2.514.20'
$ws.Range("E45").Value = '  +2.60%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '23.05'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.55%  '
$ws.Range("E47").Value = '  -0.35%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.999'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.01%  '
$ws.Range("E49").Value = '  +0.98%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.11'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.01%  '
$ws.Range("E51").Value = '  +1.58%  '
